# Scheduled runner update: refresh Market Board price data across all job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) in the Titan_Profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 733.1111
$ws.Range("J4").Value = 52
$ws.Range("L4").Value = 52
$ws.Range("N4").Value = -280
$ws.Range("H29").Value = 291.66666
$ws.Range("I29").Value = 250
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 750
$ws.Range("L29").Value = 900
$ws.Range("M29").Value = -469
$ws.Range("N29").Value = -1462
$ws.Range("H38").Value = 32.5
$ws.Range("I38").Value = 32.5
$ws.Range("K38").Value = 97.5
$ws.Range("M38").Value = 274.5
$ws.Range("H40").Value = 2275.3333
$ws.Range("I40").Value = 2210.8333
$ws.Range("J40").Value = 2533.3333
$ws.Range("K40").Value = 2210.8333
$ws.Range("L40").Value = 2533.3333
$ws.Range("M40").Value = -2035.8333
$ws.Range("N40").Value = -2883.3333
$ws.Range("H58").Value = 4528.2
$ws.Range("I58").Value = 104.44444
$ws.Range("J58").Value = 8147.636
$ws.Range("K58").Value = 313.33332
$ws.Range("L58").Value = 24442.908
$ws.Range("M58").Value = -163.33332
$ws.Range("N58").Value = -24742.908
$ws.Range("H132").Value = 18677.457
$ws.Range("I132").Value = 19307.334
$ws.Range("K132").Value = 57922.00199999999
$ws.Range("M132").Value = -55392.00199999999
$ws.Range("H135").Value = 2519.7646
$ws.Range("I135").Value = 2279.4546
$ws.Range("J135").Value = 2960.3333
$ws.Range("K135").Value = 20515.0914
$ws.Range("L135").Value = 26642.9997
$ws.Range("M135").Value = -17980.0914
$ws.Range("N135").Value = -31712.9997
$ws.Range("H141").Value = 3256.923
$ws.Range("I141").Value = 2686.25
$ws.Range("J141").Value = 10105
$ws.Range("K141").Value = 8058.75
$ws.Range("L141").Value = 30315
$ws.Range("M141").Value = -2878.75
$ws.Range("N141").Value = -40675
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3034.89
$ws.Range("I32").Value = 3034.89
$ws.Range("K32").Value = 3034.89
$ws.Range("M32").Value = -2747.89
$ws.Range("H63").Value = 13567.5
$ws.Range("I63").Value = 17335
$ws.Range("K63").Value = 17335
$ws.Range("M63").Value = -16649
$ws.Range("H64").Value = 24995
$ws.Range("J64").Value = 24995
$ws.Range("L64").Value = 24995
$ws.Range("N64").Value = -25491
$ws.Range("H66").Value = 13567.5
$ws.Range("I66").Value = 17335
$ws.Range("K66").Value = 86675
$ws.Range("M66").Value = -83243
$ws.Range("H67").Value = 24995
$ws.Range("J67").Value = 24995
$ws.Range("L67").Value = 24995
$ws.Range("N67").Value = -26711
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 5748.5264
$ws.Range("I97").Value = 6307.1763
$ws.Range("K97").Value = 6307.1763
$ws.Range("M97").Value = -5811.1763
$ws.Range("H122").Value = 2719.7856
$ws.Range("I122").Value = 2157.7
$ws.Range("J122").Value = 4125
$ws.Range("K122").Value = 6473.099999999999
$ws.Range("L122").Value = 12375
$ws.Range("M122").Value = -4023.099999999999
$ws.Range("N122").Value = -17275
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1400.15
$ws.Range("I20").Value = 1208
$ws.Range("J20").Value = 1757
$ws.Range("K20").Value = 1208
$ws.Range("L20").Value = 1757
$ws.Range("M20").Value = -961
$ws.Range("N20").Value = -2251
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H134").Value = 1722.5625
$ws.Range("I134").Value = 1381.037
$ws.Range("J134").Value = 3566.8
$ws.Range("K134").Value = 4143.111
$ws.Range("L134").Value = 10700.4
$ws.Range("M134").Value = -1608.111
$ws.Range("N134").Value = -15770.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 8000
$ws.Range("J25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("N25").Value = -8348
$ws.Range("H31").Value = 4038.3943
$ws.Range("I31").Value = 1832.7742
$ws.Range("K31").Value = 1832.7742
$ws.Range("M31").Value = -1537.7742
$ws.Range("H34").Value = 4038.3943
$ws.Range("I34").Value = 1832.7742
$ws.Range("K34").Value = 1832.7742
$ws.Range("M34").Value = -1630.7742
$ws.Range("H58").Value = 17545482
$ws.Range("I58").Value = 25642038
$ws.Range("J58").Value = 2942.6667
$ws.Range("K58").Value = 25642038
$ws.Range("L58").Value = 2942.6667
$ws.Range("M58").Value = -25641835
$ws.Range("N58").Value = -3348.6667
$ws.Range("H99").Value = 27780458
$ws.Range("I99").Value = 2593.3333
$ws.Range("J99").Value = 111114056
$ws.Range("K99").Value = 2593.3333
$ws.Range("L99").Value = 111114056
$ws.Range("M99").Value = -1095.3333
$ws.Range("N99").Value = -111117052
$ws.Range("H126").Value = 27780458
$ws.Range("I126").Value = 2593.3333
$ws.Range("J126").Value = 111114056
$ws.Range("K126").Value = 7779.999899999999
$ws.Range("L126").Value = 333342168
$ws.Range("M126").Value = -5309.999899999999
$ws.Range("N126").Value = -333347108
$ws.Range("H132").Value = 3789613
$ws.Range("I132").Value = 4387248
$ws.Range("J132").Value = 4592
$ws.Range("K132").Value = 13161744
$ws.Range("L132").Value = 13776
$ws.Range("M132").Value = -13159214
$ws.Range("N132").Value = -18836
$ws.Range("H136").Value = 17545482
$ws.Range("I136").Value = 25642038
$ws.Range("J136").Value = 2942.6667
$ws.Range("K136").Value = 76926114
$ws.Range("L136").Value = 8828.000100000001
$ws.Range("M136").Value = -76923564
$ws.Range("N136").Value = -13928.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2310.0938
$ws.Range("I5").Value = 1237.1818
$ws.Range("J5").Value = 4670.5
$ws.Range("K5").Value = 3711.5454
$ws.Range("L5").Value = 14011.5
$ws.Range("M5").Value = -3599.5454
$ws.Range("N5").Value = -14235.5
$ws.Range("H68").Value = 23286
$ws.Range("I68").Value = 25200.4
$ws.Range("J68").Value = 18500
$ws.Range("K68").Value = 75601.20000000001
$ws.Range("L68").Value = 55500
$ws.Range("M68").Value = -74790.20000000001
$ws.Range("N68").Value = -57122
$ws.Range("H71").Value = 23286
$ws.Range("I71").Value = 25200.4
$ws.Range("J71").Value = 18500
$ws.Range("K71").Value = 226803.6
$ws.Range("L71").Value = 166500
$ws.Range("M71").Value = -222747.6
$ws.Range("N71").Value = -174612
$ws.Range("H129").Value = 2946.4
$ws.Range("I129").Value = 3666.3333
$ws.Range("K129").Value = 10998.9999
$ws.Range("M129").Value = -5998.999899999999
$ws.Range("H132").Value = 1423.1428
$ws.Range("I132").Value = 1300.5
$ws.Range("J132").Value = 1472.2
$ws.Range("K132").Value = 11704.5
$ws.Range("L132").Value = 13249.8
$ws.Range("M132").Value = -9174.5
$ws.Range("N132").Value = -18309.8
$ws.Range("H135").Value = 2310.0938
$ws.Range("I135").Value = 1237.1818
$ws.Range("J135").Value = 4670.5
$ws.Range("K135").Value = 11134.6362
$ws.Range("L135").Value = 42034.5
$ws.Range("M135").Value = -8599.636200000001
$ws.Range("N135").Value = -47104.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1805.5
$ws.Range("J113").Value = 1800
$ws.Range("L113").Value = 1800
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 3000.75
$ws.Range("I122").Value = 3051
$ws.Range("J122").Value = 2850
$ws.Range("K122").Value = 9153
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -6703
$ws.Range("N122").Value = -13450
$ws.Range("H132").Value = 2530.457
$ws.Range("I132").Value = 2192.0344
$ws.Range("J132").Value = 4166.1665
$ws.Range("K132").Value = 6576.1032
$ws.Range("L132").Value = 12498.4995
$ws.Range("M132").Value = -4046.1032
$ws.Range("N132").Value = -17558.4995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3028.6
$ws.Range("I93").Value = 2857.2
$ws.Range("J93").Value = 3200
$ws.Range("K93").Value = 2857.2
$ws.Range("L93").Value = 3200
$ws.Range("M93").Value = -1609.2
$ws.Range("N93").Value = -5696
$ws.Range("H132").Value = 2429.8447
$ws.Range("I132").Value = 1682.0652
$ws.Range("J132").Value = 5296.3335
$ws.Range("K132").Value = 5046.1956
$ws.Range("L132").Value = 15889.0005
$ws.Range("M132").Value = -2516.1956
$ws.Range("N132").Value = -20949.0005
$ws.Range("H136").Value = 6162.6206
$ws.Range("I136").Value = 4665.5
$ws.Range("J136").Value = 7219.4116
$ws.Range("K136").Value = 13996.5
$ws.Range("L136").Value = 21658.2348
$ws.Range("M136").Value = -11446.5
$ws.Range("N136").Value = -26758.2348
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 718714.0600000001
$ws.Range("I81").Value = 4002600
$ws.Range("J81").Value = 4825.826
$ws.Range("K81").Value = 8005200
$ws.Range("L81").Value = 9651.652
$ws.Range("M81").Value = -8004139
$ws.Range("N81").Value = -11773.652
$ws.Range("H84").Value = 718714.0600000001
$ws.Range("I84").Value = 4002600
$ws.Range("J84").Value = 4825.826
$ws.Range("K84").Value = 40026000
$ws.Range("L84").Value = 48258.26
$ws.Range("M84").Value = -40020696
$ws.Range("N84").Value = -58866.26
$ws.Range("H122").Value = 1055.5319
$ws.Range("I122").Value = 993.2895
$ws.Range("K122").Value = 2979.8685
$ws.Range("M122").Value = -529.8685
$ws.Range("H132").Value = 2712.0322
$ws.Range("I132").Value = 2317.2126
$ws.Range("J132").Value = 3949.1333
$ws.Range("K132").Value = 6951.6378
$ws.Range("L132").Value = 11847.3999
$ws.Range("M132").Value = -4421.6378
$ws.Range("N132").Value = -16907.3999
$ws.Range("H136").Value = 2273
$ws.Range("I136").Value = 668.36365
$ws.Range("K136").Value = 2005.09095
$ws.Range("M136").Value = 544.90905